$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.299842119216919
$ws.Range("B1").Value = 3.287906169891357
$ws.Range("C1").Value = 4.665264129638672
$ws.Range("D1").Value = 2.246706247329712
$ws.Range("E1").Value = 1.501857757568359
